# CMIP5 -> CMIP6 mapping for ocean realm.
# Populate the new "cmip6-id" column (C) with the mapped CMIP6 vocabulary
# path for each CMIP5 property row, using a consistent look (Arial,
# left/center aligned, wrapped) throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCells = @("C6","C5","C8","C3","C9","C10","C12","C14","C35","C36","C38","C43")

# Build the new look on the first cell...
$first = $ws.Range($targetCells[0])
$first.Style = "Normal"
$first.Font.Name = "Arial"
$first.HorizontalAlignment = -4131
$first.VerticalAlignment = -4108
$first.WrapText = $true

# ...then fan it out to the rest via copy/paste-format so the style table
# stays minimal (one new font, one new cell format) instead of growing one
# extra format per cell.
$first.Copy()
for ($i = 1; $i -lt $targetCells.Length; $i++) {
    $ws.Range($targetCells[$i]).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Values, entered in the order the mappings were originally authored
# (controls shared-string table append order).
$ws.Range("C6").Value  = "Ocean > Key Properties > Bathymetry > Bathymetry Reference Dates"
$ws.Range("C5").Value  = "Ocean > Key Properties > Bathymetry > Ocean Bathymetry Type"
$ws.Range("C8").Value  = "Ocean > Grid > Discretisation > Horizontal > Pole Singularity Treatment"
$ws.Range("C3").Value  = "Ocean > Key Properties > General > Basic Approximations"
$ws.Range("C9").Value  = "Ocean > Grid > Discretisation > Horizontal > Scheme"
$ws.Range("C10").Value = "Ocean > Grid > Horizontal Grid > Horizontal Grid Type"
$ws.Range("C12").Value = "Ocean > Key Properties > General > Prognostic Variables"
$ws.Range("C14").Value = "Ocean > Key Properties > General > Model Family"
$ws.Range("C35").Value = "Ocean > Key Properties > Seawater Properties > Seawater Eos Type"
$ws.Range("C36").Value = "Ocean > Key Properties > Seawater Properties > Ocean Freezing Point"
$ws.Range("C38").Value = "Ocean > Key Properties > Seawater Properties > Ocean Specific Heat"
$ws.Range("C43").Value = "Ocean > Timestepping Framework > Timestepping Attributes > Time Step"

# Column widths: label column narrower, new mapping column wider to fit text.
$ws.Columns.Item(2).ColumnWidth = 42.285714285714285
$ws.Columns.Item(3).ColumnWidth = 60.285714285714285

# Leave selection where editing ended.
$ws.Range("C43").Select()
